$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 910
$ws.Cells.Item(2, 5).Value = 18
$ws.Cells.Item(2, 6).Value = 18
$ws.Cells.Item(2, 7).Value = 18
$ws.Cells.Item(2, 8).Value = 16
$ws.Cells.Item(2, 9).Value = 16
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1328
$ws.Cells.Item(2, 12).Value = 559
$ws.Cells.Item(2, 13).Value = 769
$ws.Cells.Item(2, 14).Value = 768
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 90
$ws.Cells.Item(2, 17).Value = 70
$ws.Cells.Item(2, 18).Value = -163
$ws.Cells.Item(2, 19).Value = 114
$ws.Cells.Item(2, 20).Value = 149
$ws.Cells.Item(2, 21).Value = -78
$ws.Cells.Item(2, 22).Value = 252
$ws.Cells.Item(2, 23).Value = 1.95
$ws.Cells.Item(2, 24).Value = 1.73
$ws.Cells.Item(2, 25).Value = 2.02
$ws.Cells.Item(2, 26).Value = 1.2
$ws.Cells.Item(2, 27).Value = 72.76000000000001
$ws.Cells.Item(2, 28).Value = 509.15
$ws.Cells.Item(2, 29).Value = 87
$ws.Cells.Item(2, 30).Value = 26.96
$ws.Cells.Item(2, 31).Value = 4296
$ws.Cells.Item(2, 32).Value = 0.54
$ws.Cells.Item(2, 33).Value = 10
$ws.Cells.Item(2, 34).Value = 0.43
$ws.Cells.Item(2, 35).Value = 11.46
$ws.Cells.Item(2, 36).Value = 18000000

# Row 3
$ws.Cells.Item(3, 4).Value = 940
$ws.Cells.Item(3, 5).Value = 57
$ws.Cells.Item(3, 6).Value = 57
$ws.Cells.Item(3, 7).Value = 53
$ws.Cells.Item(3, 8).Value = 43
$ws.Cells.Item(3, 9).Value = 42
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1359
$ws.Cells.Item(3, 12).Value = 573
$ws.Cells.Item(3, 13).Value = 786
$ws.Cells.Item(3, 14).Value = 785
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 90
$ws.Cells.Item(3, 17).Value = 64
$ws.Cells.Item(3, 18).Value = -22
$ws.Cells.Item(3, 19).Value = -17
$ws.Cells.Item(3, 20).Value = 10
$ws.Cells.Item(3, 21).Value = 54
$ws.Cells.Item(3, 22).Value = 238
$ws.Cells.Item(3, 23).Value = 6.01
$ws.Cells.Item(3, 24).Value = 4.53
$ws.Cells.Item(3, 25).Value = 5.45
$ws.Cells.Item(3, 26).Value = 3.17
$ws.Cells.Item(3, 27).Value = 72.90000000000001
$ws.Cells.Item(3, 28).Value = 550.36
$ws.Cells.Item(3, 29).Value = 235
$ws.Cells.Item(3, 30).Value = 11.23
$ws.Cells.Item(3, 31).Value = 4392
$ws.Cells.Item(3, 32).Value = 0.6
$ws.Cells.Item(3, 33).Value = 20
$ws.Cells.Item(3, 34).Value = 0.76
$ws.Cells.Item(3, 35).Value = 8.449999999999999
$ws.Cells.Item(3, 36).Value = 18000000

# Row 4
$ws.Cells.Item(4, 4).Value = 940
$ws.Cells.Item(4, 5).Value = 23
$ws.Cells.Item(4, 6).Value = 23
$ws.Cells.Item(4, 7).Value = 26
$ws.Cells.Item(4, 8).Value = 22
$ws.Cells.Item(4, 9).Value = 22
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 1466
$ws.Cells.Item(4, 12).Value = 651
$ws.Cells.Item(4, 13).Value = 815
$ws.Cells.Item(4, 14).Value = 813
$ws.Cells.Item(4, 15).Value = 2
$ws.Cells.Item(4, 16).Value = 90
$ws.Cells.Item(4, 17).Value = -95
$ws.Cells.Item(4, 18).Value = -20
$ws.Cells.Item(4, 19).Value = 69
$ws.Cells.Item(4, 20).Value = 12
$ws.Cells.Item(4, 21).Value = -108
$ws.Cells.Item(4, 22).Value = 310
$ws.Cells.Item(4, 23).Value = 2.41
$ws.Cells.Item(4, 24).Value = 2.31
$ws.Cells.Item(4, 25).Value = 2.71
$ws.Cells.Item(4, 26).Value = 1.54
$ws.Cells.Item(4, 27).Value = 79.88
$ws.Cells.Item(4, 28).Value = 572.66
$ws.Cells.Item(4, 29).Value = 120
$ws.Cells.Item(4, 30).Value = 36.5
$ws.Cells.Item(4, 31).Value = 4552
$ws.Cells.Item(4, 32).Value = 0.97
$ws.Cells.Item(4, 33).Value = 20
$ws.Cells.Item(4, 34).Value = 0.46
$ws.Cells.Item(4, 35).Value = 16.49
$ws.Cells.Item(4, 36).Value = 18000000

# Row 5
$ws.Cells.Item(5, 4).Value = 1084
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 15
$ws.Cells.Item(5, 7).Value = 13
$ws.Cells.Item(5, 8).Value = 17
$ws.Cells.Item(5, 9).Value = 17
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 1582
$ws.Cells.Item(5, 12).Value = 680
$ws.Cells.Item(5, 13).Value = 902
$ws.Cells.Item(5, 14).Value = 900
$ws.Cells.Item(5, 15).Value = 2
$ws.Cells.Item(5, 16).Value = 90
$ws.Cells.Item(5, 17).Value = 37
$ws.Cells.Item(5, 18).Value = -67
$ws.Cells.Item(5, 19).Value = 37
$ws.Cells.Item(5, 20).Value = 57
$ws.Cells.Item(5, 21).Value = -21
$ws.Cells.Item(5, 22).Value = 351
$ws.Cells.Item(5, 23).Value = 1.36
$ws.Cells.Item(5, 24).Value = 1.54
$ws.Cells.Item(5, 25).Value = 1.94
$ws.Cells.Item(5, 26).Value = 1.1
$ws.Cells.Item(5, 27).Value = 75.41
$ws.Cells.Item(5, 28).Value = 589.74
$ws.Cells.Item(5, 29).Value = 92
$ws.Cells.Item(5, 30).Value = 34.24
$ws.Cells.Item(5, 31).Value = 5035
$ws.Cells.Item(5, 32).Value = 0.63
$ws.Cells.Item(5, 33).Value = 20
$ws.Cells.Item(5, 34).Value = 0.63
$ws.Cells.Item(5, 35).Value = 21.48
$ws.Cells.Item(5, 36).Value = 18000000

# Row 6
$ws.Cells.Item(6, 4).Value = 838
$ws.Cells.Item(6, 5).Value = 11
$ws.Cells.Item(6, 6).Value = 11
$ws.Cells.Item(6, 7).Value = -17
$ws.Cells.Item(6, 8).Value = -8
$ws.Cells.Item(6, 9).Value = -8
$ws.Cells.Item(6, 11).Value = 1491
$ws.Cells.Item(6, 12).Value = 653
$ws.Cells.Item(6, 13).Value = 838
$ws.Cells.Item(6, 14).Value = 837
$ws.Cells.Item(6, 16).Value = 90
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = -35
$ws.Cells.Item(6, 19).Value = 28
$ws.Cells.Item(6, 20).Value = 10
$ws.Cells.Item(6, 21).Value = 14
$ws.Cells.Item(6, 22).Value = 383
$ws.Cells.Item(6, 23).Value = 1.35
$ws.Cells.Item(6, 24).Value = -0.98
$ws.Cells.Item(6, 25).Value = -0.9399999999999999
$ws.Cells.Item(6, 26).Value = -0.53
$ws.Cells.Item(6, 27).Value = 77.84999999999999
$ws.Cells.Item(6, 28).Value = 635.76
$ws.Cells.Item(6, 29).Value = -45
$ws.Cells.Item(6, 30).Value = -97.95
$ws.Cells.Item(6, 31).Value = 4682
$ws.Cells.Item(6, 32).Value = 0.9399999999999999
$ws.Cells.Item(6, 33).Value = 20
$ws.Cells.Item(6, 34).Value = 0.45
$ws.Cells.Item(6, 35).Value = -44.05
$ws.Cells.Item(6, 36).Value = 18000000

# Clear forecast rows 7-9 (columns D:AJ), keep only A/B/C
$ws.Range("D7:AJ9").ClearContents()
